# "connect all prodigi & proaktif"
#
# Adds four more report sheets (Sheet2..Sheet5), each a copy of the
# existing "Sheet1" header template (No / Singer / Title / Total trafic /
# total revenue / Pencipta / Partner / Artis / Nama Chanel Marketing /
# Revenue Prodigi), and leaves the last one ("Sheet5") as the active tab -
# mirroring how BA/Sheet1 used to be wired up.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item("Sheet1")

# Column widths (in "characters") used by the template header row. COM's
# ColumnWidth setter is expressed relative to the ~0.8333-char padding
# Excel always adds on top, so back that out before assigning.
$pad = 0.8333333333333334
$colWidths = @(3.5703125, 6.5703125, 5, 10.42578125, 13.140625, 8.7109375, 7.5703125, 5.140625, 22.5703125, 15.7109375)

$headers = @("No", "Singer", "Title", "Total trafic", "total revenue", "Pencipta", "Partner", "Artis", "Nama Chanel Marketing", "Revenue Prodigi")

function New-ReportSheet([string]$name) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add($null, $lastSheet)
    $ws.Name = $name

    # Bring over the header row's formatting (style) from Sheet1, then
    # (re)write the literal header text on top of it.
    [void]$sheet1.Range("A1:J1").Copy()
    [void]$ws.Range("A1:J1").PasteSpecial(-4122)   # xlPasteFormats

    for ($i = 0; $i -lt $headers.Count; $i++) {
        $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
    }

    for ($i = 0; $i -lt $colWidths.Count; $i++) {
        $ws.Columns.Item($i + 1).ColumnWidth = $colWidths[$i] - $pad
    }

    return $ws
}

$sheet2 = New-ReportSheet "Sheet2"
[void]$sheet2.Cells.Select()

$sheet3 = New-ReportSheet "Sheet3"
[void]$sheet3.Cells.Select()

$sheet4 = New-ReportSheet "Sheet4"
[void]$sheet4.Cells.Select()

$sheet5 = New-ReportSheet "Sheet5"
[void]$sheet5.Range("E8").Select()

# Sheet1's own selection moves from a single cell to the whole header row,
# and it picks up a (portrait, paper size 9) page setup.
$sheet1.PageSetup.PaperSize = 9
$sheet1.PageSetup.Orientation = 1
[void]$sheet1.Rows.Item(1).Select()

# Sheet5 (the newest sheet, 6th tab) becomes the active tab; this also
# clears tabSelected from BA (which used to carry it).
[void]$sheet5.Activate()
